$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.082.05"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "1.749.42"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "'233.40"
$ws.Range("E5").Value = "  +3.44%  "
$ws.Range("D6").Value = "'0.9983"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("D7").Value = "'0.5246"
$ws.Range("E7").Value = "  +1.95%  "
$ws.Range("D8").Value = "'0.2782"
$ws.Range("E8").Value = "  +3.98%  "
$ws.Range("D9").Value = "'40.17"
$ws.Range("E9").Value = "  +1.80%  "
$ws.Range("D10").Value = "'0.06193"
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("D11").Value = "1.752.40"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").Value = "'0.07168"
$ws.Range("E12").Value = "  +2.76%  "
$ws.Range("D13").Value = "'15.35"
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("D14").Value = "'0.6452"
$ws.Range("E14").Value = "  +3.15%  "
$ws.Range("D15").Value = "'4.569"
$ws.Range("E15").Value = "  +2.35%  "
$ws.Range("D16").Value = "'78.15"
$ws.Range("E16").Value = "  +2.75%  "
$ws.Range("D17").Value = "'0.9984"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "'0.9984"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").Value = "25.979.04"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("E20").Value = "  +2.57%  "
$ws.Range("D21").Value = "'0.000006686"
$ws.Range("E21").Value = "  +3.23%  "
$ws.Range("D22").Value = "1.973.19"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'4.301"
$ws.Range("E23").Value = "  +7.65%  "
$ws.Range("D24").Value = "'8.808"
$ws.Range("E24").Value = "  +5.45%  "
$ws.Range("D25").Value = "'5.211"
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("D26").Value = "'138.84"
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").Value = "'15.34"
$ws.Range("E28").Value = "  +3.09%  "
$ws.Range("D29").Value = "'1.817"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "'104.21"
$ws.Range("E30").Value = "  +2.06%  "
$ws.Range("D31").Value = "'0.08341"
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("D32").Value = "'3.771"
$ws.Range("E32").Value = "  +4.46%  "
$ws.Range("D33").Value = "'3.625"
$ws.Range("E33").Value = "  +8.21%  "
$ws.Range("E34").Value = "  +3.54%  "
$ws.Range("D35").Value = "'2.622"
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("D36").Value = "'1.004"
$ws.Range("E36").Value = "  +3.78%  "
$ws.Range("D37").Value = "'0.6306"
$ws.Range("E37").Value = "  +7.02%  "
$ws.Range("D38").Value = "'2.705"
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("D39").Value = "'0.01599"
$ws.Range("E39").Value = "  +2.54%  "
$ws.Range("D40").Value = "'1.941"
$ws.Range("E40").Value = "  +2.36%  "
$ws.Range("D41").Value = "'0.9976"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").Value = "'98.52"
$ws.Range("E42").Value = "  -2.95%  "
$ws.Range("D43").Value = "'0.3912"
$ws.Range("E43").Value = "  +3.93%  "
$ws.Range("D44").Value = "'0.7366"
$ws.Range("E44").Value = "  +1.56%  "
$ws.Range("D45").Value = "'5.069"
$ws.Range("E45").Value = "  +5.57%  "
$ws.Range("D46").Value = "'0.1134"
$ws.Range("E46").Value = "  +3.81%  "
$ws.Range("D47").Value = "'6.312"
$ws.Range("E47").Value = "  +2.58%  "
$ws.Range("D48").Value = "'0.05357"
$ws.Range("E48").Value = "  -2.10%  "
$ws.Range("D49").Value = "'54.01"
$ws.Range("E49").Value = "  +4.97%  "
$ws.Range("D50").Value = "'30.41"
$ws.Range("E50").Value = "  +2.91%  "
$ws.Range("D51").Value = "'7.653"
$ws.Range("E51").Value = "  +4.17%  "
